$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) column values are stored as text in the source sheet (e.g. "1.00", "11.50").
# Force text number format on those specific cells before assignment so Excel does not
# auto-coerce them to numbers (which would drop significant trailing zeros, e.g. 11.50 -> 11.5).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.144.03"
$ws.Range("E2").Value = "  +2.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.529.78"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.81"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.06"
$ws.Range("E6").Value = "  +2.39%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.529.12"
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("E11").Value = "  +2.86%  "
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.347"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.98"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.991.24"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.014.26"
$ws.Range("E17").Value = "  +2.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.531.93"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.01"
$ws.Range("E19").Value = "  +3.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.50"
$ws.Range("E20").Value = "  +2.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "364.20"
$ws.Range("E21").Value = "  +4.43%  "
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.70"
$ws.Range("E23").Value = "  +2.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.95"
$ws.Range("E24").Value = "  -0.88%  "
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.00"
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.29"
$ws.Range("E27").Value = "  +3.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.660.19"
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0993"
$ws.Range("E30").Value = "  +1.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "547.25"
$ws.Range("E31").Value = "  +3.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.27"
$ws.Range("E32").Value = "  +1.93%  "
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("E34").Value = "  +2.45%  "
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.62"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.89"
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("E40").Value = "  +1.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.357"
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.21"
$ws.Range("E42").Value = "  +2.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.80"
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.563"
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "147.45"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("E48").Value = "  +1.92%  "
$ws.Range("E49").Value = "  +2.46%  "
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0759"
$ws.Range("E51").Value = "  +1.12%  "
